$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.046.01"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.311.42"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.97"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.65"
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.522"
$ws.Range("E7").Value = "  +3.48%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.67"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.93"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.90"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "2.672.72"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "2.264.33"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("E17").Value = "  -3.17%  "
$ws.Range("D18").Value = "42.969.95"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.53"
$ws.Range("E19").Value = "  +7.27%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.34"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.88"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.15"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.18"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("E30").Value = "  -12.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.35"
$ws.Range("E31").Value = "  -3.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.20"
$ws.Range("E32").Value = "  +3.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +2.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.27"
$ws.Range("E35").Value = "  +6.45%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.80"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.112"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").Value = "2.000.13"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0288"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.16"
$ws.Range("E44").Value = "  -6.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.11"
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.60"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.83"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.88"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("D49").Value = "2.538.81"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.99"
$ws.Range("E50").Value = "  +5.15%  "
$ws.Range("E51").Value = "  +0.18%  "
